# Measurements.xlsx - add two new measurement cycles (rows 31 & 32) and
# fill in the previously-blank row 30, extend Table1 / conditional
# formatting / dimension accordingly, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 30 - was just empty styled placeholder cells, now real data.
# ---------------------------------------------------------------------
$ws.Range("B30").Value = 27
$ws.Range("E30").Value = 1408.7
$ws.Range("F30").Value = 1407.33
$ws.Range("G30").Value = 45894.666666666664
$ws.Range("H30").Value = 45895.46597222222
$ws.Range("I30").Value = 9.6999999999999993
$ws.Range("J30").Value = 8.5500000000000007
$ws.Range("K30").Value = 2.5

# ---------------------------------------------------------------------
# 2. Row 31 - brand-new row; copy formats from row 30 for the cells
#    that need a non-default style (date cells + formula cells), then
#    fill in the values.
# ---------------------------------------------------------------------
$ws.Range("G30").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("H30").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("O30").Copy()
$ws.Range("O31").PasteSpecial(-4122)
$ws.Range("P30").Copy()
$ws.Range("P31").PasteSpecial(-4122)
$ws.Range("Q30").Copy()
$ws.Range("Q31").PasteSpecial(-4122)
$ws.Range("R30").Copy()
$ws.Range("R31").PasteSpecial(-4122)
$ws.Range("S30").Copy()
$ws.Range("S31").PasteSpecial(-4122)
$ws.Range("U30").Copy()
$ws.Range("U31").PasteSpecial(-4122)
$ws.Range("V30").Copy()
$ws.Range("V31").PasteSpecial(-4122)

$ws.Range("B31").Value = 28
$ws.Range("C31").Value = 286
$ws.Range("D31").Value = 159
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 45894.665277777778
$ws.Range("H31").Value = 45895.466666666667
$ws.Range("I31").Value = 11.5
$ws.Range("J31").Value = 9.9
$ws.Range("K31").Value = 2.5

# ---------------------------------------------------------------------
# 3. Row 32 - another brand-new row, same treatment.
# ---------------------------------------------------------------------
$ws.Range("G30").Copy()
$ws.Range("G32").PasteSpecial(-4122)
$ws.Range("H30").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("O30").Copy()
$ws.Range("O32").PasteSpecial(-4122)
$ws.Range("P30").Copy()
$ws.Range("P32").PasteSpecial(-4122)
$ws.Range("Q30").Copy()
$ws.Range("Q32").PasteSpecial(-4122)
$ws.Range("R30").Copy()
$ws.Range("R32").PasteSpecial(-4122)
$ws.Range("S30").Copy()
$ws.Range("S32").PasteSpecial(-4122)
$ws.Range("U30").Copy()
$ws.Range("U32").PasteSpecial(-4122)
$ws.Range("V30").Copy()
$ws.Range("V32").PasteSpecial(-4122)

$ws.Range("B32").Value = 29
$ws.Range("C32").Value = 286
$ws.Range("D32").Value = 159
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 45895.534722222219
$ws.Range("H32").Value = 45895.573611111111
$ws.Range("I32").Value = 9.85
$ws.Range("J32").Value = 8.9
$ws.Range("K32").Value = 25

# ---------------------------------------------------------------------
# 4. Formula columns - O, P, Q, R, S, V get shared formulas spanning
#    the two fully-new rows (30:31 continues the existing O30 stub,
#    32 is its own single-cell shared formula, matching how Excel
#    would fill these down one block at a time); U is always a
#    standalone (non-shared) formula with a per-row literal constant.
# ---------------------------------------------------------------------
$ws.Range("O30:O31").Formula = "=H30-G30"
$ws.Range("P30:P31").Formula = "=O30"
$ws.Range("Q30:Q31").Formula = "=I30-J30"
$ws.Range("R30:R31").Formula = "=(F30-E30)/0.9982"
$ws.Range("S30:S31").Formula = "=K30*P30"
$ws.Range("U30").Formula = "=Q30*1440/1151"
$ws.Range("V30:V31").Formula = "=(1-ABS(U30-K30)/K30)*100%"

$ws.Range("U31").Formula = "=Q31*1440/1154"

$ws.Range("O32").Formula = "=H32-G32"
$ws.Range("P32").Formula = "=O32"
$ws.Range("Q32").Formula = "=I32-J32"
$ws.Range("R32").Formula = "=(F32-E32)/0.9982"
$ws.Range("S32").Formula = "=K32*P32"
$ws.Range("U32").Formula = "=Q32*1440/56"
$ws.Range("V32").Formula = "=(1-ABS(U32-K32)/K32)*100%"

# ---------------------------------------------------------------------
# 5. Extend Table1 to cover the two new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("B3:L32"))

# ---------------------------------------------------------------------
# 6. Extend the conditional-formatting ranges that used to stop at
#    row 30 so they now cover row 32 as well.
# ---------------------------------------------------------------------
$ws.Range("V4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("V4:V32"))
$ws.Range("P4").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("P4:P32"))

# ---------------------------------------------------------------------
# 7. Move the active selection.
# ---------------------------------------------------------------------
$ws.Range("P26").Select()
